$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row of data to the registration table (row 3)
$ws.Range("A3").Value = "The Owls"
$ws.Range("C3").Value = "Quản lý khách hàng(xem, sửa, thêm)"
$ws.Range("D3").Value = "Ngô Thị Mai Lý (1412310)"
$ws.Range("B3").Value = "21/12/2017"

# Update the active selection to the newly added cell, as in the saved file
$ws.Range("B3").Select()
